$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Marion Queen Ramos"

# Row 7 - __init__ / Attribute set to input values (valid inputs)
$ws.Range("E7").Value = "None, just creating a new object."
$ws.Range("F7").Value = "A title, author, and genre that are all valid (e.g., title=""1985"", author=""Orson Scott Card"", genre=Genre.FICTION)."
$ws.Range("G7").Value = "The object gets created just fine with those values."

# Row 8 - Exception raised when title is blank
$ws.Range("E8").Value = "None, just trying to create a new object."
$ws.Range("F8").Value = "title="""", author=""Orson Scott Card"", genre=Genre.FICTION"
$ws.Range("G8").Value = "It should throw a ValueError saying ""Title cannot be blank."""

# Row 9 - Exception raised when author is blank
$ws.Range("E9").Value = "None, just trying to create a new object."
$ws.Range("F9").Value = " title=""1985"", author="""", genre=Genre.FICTION"
$ws.Range("G9").Value = "It should throw a ValueError saying ""Author cannot be blank."""

# Row 10 - Exception raised when invalid Genre
$ws.Range("E10").Value = "None, just trying to create a new object."
$ws.Range("F10").Value = " title=""1985"", author=""Orson Scott Card"", genre=""INVALID_GENRE"" (or YOUNG ADULT)"
$ws.Range("G10").Value = "It should throw a ValueError saying ""Invalid Genre."""

# Row 11 - returns title attribute
$ws.Range("E11").Value = "The object should already exist"
$ws.Range("F11").Value = "None, just calling the title property."
$ws.Range("G11").Value = "Should return the title, like ""1985""."

# Row 12 - returns author attribute
$ws.Range("E12").Value = "The object should already exist."
$ws.Range("F12").Value = "None, just calling the author property."
$ws.Range("G12").Value = "Should return the author, like ""Orson Scott Card""."

# Row 13 - returns Genre attribute
$ws.Range("E13").Value = " The object should already exist."
$ws.Range("F13").Value = "None, just calling the genre property"
$ws.Range("G13").Value = "Should return the genre, like Genre.FICTION."
